$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26; this shifts existing rows 26-33 down to 27-34
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
# Values mirror the other "Locoto" records in this sheet/column layout.
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value2 = 44776
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = 100112042
$ws.Range("G26").Value = "Locoto"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 2700
$ws.Range("L26").Value = 2700
$ws.Range("M26").Value = 2700
$ws.Range("N26").Value = "$/kilo"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 2700
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"
